$wb = $excel.ActiveWorkbook

# The workbook currently ends with a "总计" (grand-total) summary sheet.
# We insert a brand-new "2022-Q1" sheet just in front of it (so the sheet
# order stays chronological) and then refresh the summary sheet so it
# includes the new quarter too.

$totalSheetNameHandle = $wb.Worksheets.Item($wb.Worksheets.Count)

# --- 1. Create the new "2022-Q1" detail sheet -----------------------------
$q1 = $wb.Worksheets.Add($totalSheetNameHandle)
$q1.Name = "2022-Q1"

# NOTE: worksheet handles resolve by position, so the handle obtained above
# (pointing at the last tab, "总计") now refers to the freshly-inserted
# "2022-Q1" sheet instead, since it took over that position. Re-fetch the
# "总计" sheet by name so later writes land on the right tab.
$totalSheet = $wb.Worksheets.Item("总计")

# Match the page layout (margins) and outline summary direction used by
# all the other sheets in this workbook.
$q1.PageSetup.LeftMargin = 54
$q1.PageSetup.RightMargin = 54
$q1.PageSetup.TopMargin = 72
$q1.PageSetup.BottomMargin = 72
$q1.PageSetup.HeaderMargin = 36
$q1.PageSetup.FooterMargin = 36
$q1.Outline.SummaryRow = 1
$q1.Outline.SummaryColumn = 1

# Header row (row 1) - column A has no header, matching the other quarter sheets.
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Apply the bold/centered/bordered header style (copied from an existing
# quarter sheet's header row) to row 1 of the new sheet.
$wb.Worksheets.Item(1).Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

# Columns B (fund code) and D-G (scale / position figures) must stay text
# even though they look numeric, so values such as the leading zero in
# "010874" or the trailing zero in "0.10" survive. Mark those ranges as
# text *before* assigning values so Excel doesn't auto-convert them to
# numbers; column C (fund name) is never numeric-looking so needs no
# special handling, and column H (rank) is a genuine number.
$q1.Range("B2:B15").NumberFormat = "@"
$q1.Range("D2:G15").NumberFormat = "@"

# Data rows (row index, A value, B..H values). Columns B-G hold text
# (fund code / name / figures retained as strings so trailing zeros such
# as "0.10" are preserved); column H holds a numeric rank.
$q1Rows = @(
    ,@(2,  0,  "010874", "泰康品质生活混合A",             "13.17", "81.43", "3.21", "0.4228", 5)
    ,@(3,  1,  "005014", "泰康景泰回报混合A",             "11.64", "27.21", "1.52", "0.1769", 4)
    ,@(4,  2,  "010965", "中银鑫新消费成长混合A",         "4.76",  "73.66", "3.57", "0.1699", 4)
    ,@(5,  3,  "010875", "泰康品质生活混合C",             "4.39",  "81.43", "3.21", "0.1409", 5)
    ,@(6,  4,  "009414", "中银大健康股票A",               "2.96",  "86.25", "3.36", "0.0995", 7)
    ,@(7,  5,  "012071", "中加喜利回报一年持有期混合A",   "5.20",  "46.57", "1.68", "0.0874", 10)
    ,@(8,  6,  "005775", "中加转型动力灵活配置混合A",     "3.41",  "66.34", "2.43", "0.0829", 6)
    ,@(9,  7,  "009242", "中加核心智造混合A",             "2.05",  "65.71", "3.35", "0.0687", 5)
    ,@(10, 8,  "005776", "中加转型动力灵活配置混合C",     "1.92",  "66.34", "2.43", "0.0467", 6)
    ,@(11, 9,  "012072", "中加喜利回报一年持有期混合C",   "2.74",  "46.57", "1.68", "0.0460", 10)
    ,@(12, 10, "010962", "中银鑫新消费成长混合C",         "0.82",  "73.66", "3.57", "0.0293", 4)
    ,@(13, 11, "005015", "泰康景泰回报混合C",             "0.63",  "27.21", "1.52", "0.0096", 4)
    ,@(14, 12, "009243", "中加核心智造混合C",             "0.10",  "65.71", "3.35", "0.0034", 5)
    ,@(15, 13, "010321", "中银大健康股票C",               "0.10",  "86.25", "3.36", "0.0034", 7)
)

foreach ($r in $q1Rows) {
    $row = $r[0]
    $q1.Cells.Item($row, 1).Value = $r[1]
    $q1.Cells.Item($row, 2).Value = $r[2]
    $q1.Cells.Item($row, 3).Value = $r[3]
    $q1.Cells.Item($row, 4).Value = $r[4]
    $q1.Cells.Item($row, 5).Value = $r[5]
    $q1.Cells.Item($row, 6).Value = $r[6]
    $q1.Cells.Item($row, 7).Value = $r[7]
    $q1.Cells.Item($row, 8).Value = $r[8]
}

# Give column A of the data rows the same style used in row 1 / other sheets
# (bold, centered, bordered).
$wb.Worksheets.Item(1).Range("A2").Copy()
$q1.Range("A2:A15").PasteSpecial(-4122)

# --- 2. Refresh the "总计" summary sheet with the new quarter -------------
# A new row is inserted at the top of the data (row 2) for "2022-Q1"; the
# existing quarters shift down by one row and their running index (column A)
# is renumbered 0..5.
$totalRows = @(
    ,@(2, 0, "2022-Q1", 14, 1.39)
    ,@(3, 1, "2021-Q4", 15, 1.26)
    ,@(4, 2, "2021-Q3", 9,  0.96)
    ,@(5, 3, "2021-Q2", 6,  0.58)
    ,@(6, 4, "2021-Q1", 18, 1.34)
    ,@(7, 5, "2020-Q4", 8,  0.31)
)

foreach ($r in $totalRows) {
    $row = $r[0]
    $totalSheet.Cells.Item($row, 1).Value = $r[1]
    $totalSheet.Cells.Item($row, 2).Value = $r[2]
    $totalSheet.Cells.Item($row, 3).Value = $r[3]
    $totalSheet.Cells.Item($row, 4).Value = $r[4]
}

# Row 7 is brand new on this sheet - give its column-A cell the same style
# as the rest of column A (bold, centered, bordered).
$totalSheet.Range("A6").Copy()
$totalSheet.Range("A7").PasteSpecial(-4122)
